$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates (rich-text shared strings collapse to plain text
# when edited through this object model, but the cell style already
# carries the same font/size that the original runs used, so the
# rendered appearance is unchanged).
# ---------------------------------------------------------------------
$ws.Range("A8").Value2 = "Volume 31   Number  28"
$ws.Range("C9").Value2 = "Report Covering the Week  7/8/2024  Through  7/14/2024"

# ---------------------------------------------------------------------
# Helper: a couple of cells flip between a numeric style and the
# plain-text "placeholder" style (s=14, shared strings "0"/"***.*").
# Copying the already-correctly-styled neighbour cell over first makes
# the destination pick up both the right style index AND (for the
# string case) the right shared-string value; a subsequent .Value
# write (when needed) only changes the content, not the style.
# ---------------------------------------------------------------------

# Row 15 - Rape
$ws.Range("D15").Copy($ws.Range("C15"))   # C15 -> "0" (string placeholder)
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100

# Row 16 - Robbery
$ws.Range("D15").Copy($ws.Range("C16"))   # C16 -> "0" (string placeholder)
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -75
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = -24
$ws.Range("L16").Value = -15.555555555555
$ws.Range("M16").Value = -2.564102564102
$ws.Range("N16").Value = -88.343558282208

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("I17").Value = 63
$ws.Range("J17").Value = 53
$ws.Range("K17").Value = 18.867924528301
$ws.Range("L17").Value = 8.620689655172
$ws.Range("M17").Value = 75
$ws.Range("N17").Value = -33.684210526315

# Row 18 - Burglary
$ws.Range("C18").Value = 1
$ws.Range("C18").Copy($ws.Range("D18"))   # D18 -> numeric style (s=15)
$ws.Range("D18").Value = 3
$ws.Range("H18").Copy($ws.Range("E18"))   # E18 -> percent style (s=16)
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 45
$ws.Range("J18").Value = 62
$ws.Range("K18").Value = -27.419354838709
$ws.Range("L18").Value = -60.176991150442
$ws.Range("M18").Value = -22.413793103448
$ws.Range("N18").Value = -92.512479201331

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 35
$ws.Range("H19").Value = -27.083333333333
$ws.Range("I19").Value = 303
$ws.Range("J19").Value = 390
$ws.Range("K19").Value = -22.307692307692
$ws.Range("L19").Value = -9.552238805970
$ws.Range("M19").Value = -18.766756032171
$ws.Range("N19").Value = -73.743500866551

# Row 20 - G.L.A.
$ws.Range("D20").Value = 3
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -83.333333333333
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = -56.666666666666
$ws.Range("L20").Value = -58.064516129032
$ws.Range("N20").Value = -96.209912536443

# Row 21 - TOTAL
$ws.Range("F21").Value = 53
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = -31.168831168831
$ws.Range("I21").Value = 469
$ws.Range("J21").Value = 590
$ws.Range("K21").Value = -20.508474576271
$ws.Range("L21").Value = -20.238095238095
$ws.Range("M21").Value = -10.325047801147
$ws.Range("N21").Value = -81.433095803642

# Row 22 - Transit
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("L22").Value = 71.428571428571
$ws.Range("M22").Value = 26.315789473684

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = -51.219512195122
$ws.Range("F24").Value = 82
$ws.Range("G24").Value = 135
$ws.Range("H24").Value = -39.259259259259
$ws.Range("I24").Value = 517
$ws.Range("J24").Value = 567
$ws.Range("K24").Value = -8.818342151675
$ws.Range("L24").Value = -23.633677991137
$ws.Range("M24").Value = 52.507374631268

# Row 25 - Retail Theft
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 36
$ws.Range("E25").Value = -58.333333333333
$ws.Range("F25").Value = 62
$ws.Range("G25").Value = 117
$ws.Range("H25").Value = -47.008547008547
$ws.Range("I25").Value = 419
$ws.Range("J25").Value = 469
$ws.Range("K25").Value = -10.660980810234
$ws.Range("L25").Value = -25.577264653641

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 12.5
$ws.Range("I26").Value = 109
$ws.Range("J26").Value = 142
$ws.Range("K26").Value = -23.239436619718
$ws.Range("L26").Value = -11.382113821138
$ws.Range("M26").Value = -16.793893129771

# Row 27 - UCR Rape*
$ws.Range("D27").Copy($ws.Range("C27"))   # C27 -> "0" (string placeholder)
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100

# Row 28 - Other Sex Crimes
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 2
$ws.Range("J28").Value = 32
$ws.Range("K28").Value = -3.125

# Row 31 - Hate Crimes
$ws.Range("F31").Value = 3
$ws.Range("H31").Value = 50
$ws.Range("I31").Value = 12
$ws.Range("K31").Value = 140
$ws.Range("L31").Value = 100
